# Update column G ("K") values on Sheet1, rows 2-25, to reflect the
# regenerated save_data (author's commit: "regen save_data to use K
# instead of Strike#, regen std/mean, calc and write s_vals").
#
# Only the values in column G change; every other cell is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 2
    4  = 0
    5  = 3
    7  = 0
    8  = 3
    9  = 1
    10 = 2
    11 = 1
    12 = 2
    13 = 2
    14 = 2
    15 = 1
    16 = 0
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 1
    22 = 3
    23 = 0
    24 = 1
    25 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
